# Updates the cryptos list with new prices / volume percentages, and
# shifts the last rows (BabyDogeCoin dropped, Aave/Cronos/Mantle/EnergySwap
# move up one row, USDD added at the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple D/E (Price / Volume(1h)) value updates for rows that keep their coin ---
$updates = @(
    @{ Row = 2;  D = "26.504.16";  E = "  -0.65%  " },
    @{ Row = 3;  D = "1.627.29";   E = "  -0.45%  " },
    @{ Row = 4;  D = $null;        E = "  +0.12%  " },
    @{ Row = 5;  D = "212.82";     E = "  +0.03%  " },
    @{ Row = 6;  D = $null;        E = "  +1.23%  " },
    @{ Row = 7;  D = $null;        E = "  +0.10%  " },
    @{ Row = 8;  D = "0.0623";     E = "  -0.03%  " },
    @{ Row = 9;  D = "0.248";      E = "  -1.60%  " },
    @{ Row = 10; D = "18.78";      E = "  -1.35%  " },
    @{ Row = 11; D = "0.0843";     E = "  +0.55%  " },
    @{ Row = 12; D = "1.851.68";   E = "  -0.59%  " },
    @{ Row = 13; D = "1.667.90";   E = "  +2.00%  " },
    @{ Row = 14; D = $null;        E = "  +1.58%  " },
    @{ Row = 15; D = "0.523";      E = "  -0.48%  " },
    @{ Row = 16; D = "64.97";      E = "  +3.15%  " },
    @{ Row = 17; D = "26.513.61";  E = "  -0.57%  " },
    @{ Row = 18; D = $null;        E = "  +0.00%  " },
    @{ Row = 19; D = "214.29";     E = "  +2.83%  " },
    @{ Row = 20; D = $null;        E = "  +0.14%  " },
    @{ Row = 21; D = "4.29";       E = "  -0.50%  " },
    @{ Row = 22; D = "6.26";       E = "  +1.55%  " },
    @{ Row = 23; D = "9.28";       E = "  -1.18%  " },
    @{ Row = 24; D = "2.08";       E = "  +8.57%  " },
    @{ Row = 25; D = "148.67";     E = "  +1.52%  " },
    @{ Row = 26; D = "1.00";       E = "  +0.14%  " },
    @{ Row = 27; D = "0.120";      E = "  -0.30%  " },
    @{ Row = 28; D = "6.86";       E = "  +1.88%  " },
    @{ Row = 29; D = "15.51";      E = "  +0.86%  " },
    @{ Row = 30; D = $null;        E = "  -1.59%  " },
    @{ Row = 31; D = $null;        E = "  -0.96%  " },
    @{ Row = 32; D = $null;        E = "  +3.04%  " },
    @{ Row = 33; D = $null;        E = "  -0.40%  " },
    @{ Row = 34; D = "1.237.78";   E = "  +5.91%  " },
    @{ Row = 35; D = $null;        E = "  +0.05%  " },
    @{ Row = 36; D = $null;        E = "  -1.93%  " },
    @{ Row = 37; D = "0.0174";     E = "  +4.25%  " },
    @{ Row = 38; D = $null;        E = "  +0.10%  " },
    @{ Row = 39; D = "0.508";      E = "  +0.92%  " },
    @{ Row = 40; D = "0.794";      E = "  -1.43%  " },
    @{ Row = 41; D = $null;        E = "  -1.93%  " },
    @{ Row = 42; D = "0.796";      E = "  +0.27%  " },
    @{ Row = 43; D = $null;        E = "  -0.66%  " },
    @{ Row = 44; D = "1.761.31";   E = "  -0.68%  " },
    @{ Row = 45; D = "92.93";      E = "  +0.62%  " },
    @{ Row = 46; D = "1.59";       E = "  +2.42%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# --- Rows 47-51: BabyDogeCoin is removed from the ranking, so every coin
#     previously below it shifts up one row, and USDD newly enters at the
#     bottom (row 51). Overwrite columns B (Coin), C (Link), D (Price) and
#     E (Volume(1h)) for these rows with the new lineup. ---
$lastRows = @(
    @{ Row = 47; B = "Aave";        C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";               D = "54.90";  E = "  +0.36%  " },
    @{ Row = 48; B = "Cronos";      C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";           D = "0.0509"; E = "  -0.64%  " },
    @{ Row = 49; B = "Mantle";      C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt";              D = "0.406";  E = "  -1.02%  " },
    @{ Row = 50; B = "EnergySwap";  C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";          D = "7.46";   E = "  -0.75%  " },
    @{ Row = 51; B = "USDD";        C = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd";               D = "1.00";   E = "  +0.09%  " }
)

foreach ($r in $lastRows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $dCell = $ws.Cells.Item($r.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
